$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-0.29***"
$ws.Range("B3").Value = "-1.02*"
$ws.Range("B4").Value = "0.06***"
$ws.Range("B5").Value = "0.01***"

$ws.Range("C2").Value = "-0.02***"
$ws.Range("C3").Value = "-0.42***"
$ws.Range("C4").Value = "-0.0*"
$ws.Range("C5").Value = "0.0***"

$ws.Range("D2").Value = "0.26*"
$ws.Range("D4").Value = "0.32***"
$ws.Range("D5").Value = "-0.01*"

$ws.Range("E2").Value = "-3.7*"

# These values have no "*" suffix, so Excel would otherwise read them as
# numbers. Force text storage (matching the source's shared-string text
# cells) by switching to a text number format just long enough to type
# them in, then restoring the cell's normal (default) style.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.65"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.49"
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.53"
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.1"
$ws.Range("E5").Style = "Normal"
